# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 (R) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 187
$wsOff.Range("C3").Value = 128
$wsOff.Range("D3").Value = 44
$wsOff.Range("E3").Value = 18
$wsOff.Range("F3").Value = 7

# DEF sheet - row 3 (R) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 176
$wsDef.Range("C3").Value = 115
$wsDef.Range("D3").Value = 54
$wsDef.Range("E3").Value = 29
